# Insert a new data row at row 87 (pushing existing rows 87..151 down to 88..152)
# and populate it with the new price observation, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 87..151 down to 88..152
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with its data
$ws.Cells.Item(87, 1).Value = 8
$ws.Cells.Item(87, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 44762
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = 100112001
$ws.Cells.Item(87, 7).Value = "Berenjena"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 460
$ws.Cells.Item(87, 11).Value = 10000
$ws.Cells.Item(87, 12).Value = 11000
$ws.Cells.Item(87, 13).Value = 10500
$ws.Cells.Item(87, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(87, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(87, 16).Value = 210
$ws.Cells.Item(87, 17).Value = 50
$ws.Cells.Item(87, 18).Value = "Hortaliza"
